$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# Row 10 - B/C text replaced (Objetivos body -> professor name placeholder)
# B10/C10 already exist, plain Value assignment keeps the correct column style.
$ws.Range("B10").Value = "11079086 - Herlandí de Souza Andrade"
$ws.Range("C10").Value = "11079086 - Herlandí de Souza Andrade"

# Row 13 - now carries "Programa resumido:" / "Semestral", becomes a 60pt row.
# A13/B13/C13 are brand-new cells for this row.
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B10").Copy()
$ws.Range("B13").PasteSpecial($xlPasteFormats)
$ws.Range("B13").Value = "Semestral"
$ws.Range("C10").Copy()
$ws.Range("C13").PasteSpecial($xlPasteFormats)
$ws.Range("C13").Value = "Semestral"
$ws.Rows(13).RowHeight = 60

# Row 14 - now "Short syllabus:" content (cells already existed)
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Value = "1. Management Practice Areas. 2. Organizational structure"
$ws.Range("C14").Value = "1. Management Practice Areas. 2. Organizational structure"

# Row 15 - now "Programa:" / activation date, becomes 120pt row.
# The date-like text must be copied from B8/C8 (plain Paste, not a typed Value
# assignment) so the engine doesn't reinterpret "01/01/2021" as a date serial.
$ws.Range("A15").Value = "Programa:"
$ws.Range("B8").Copy()
$ws.Paste($ws.Range("B15"))
$ws.Range("C8").Copy()
$ws.Paste($ws.Range("C15"))
$ws.Rows(15).RowHeight = 120

# Row 16 - now "Syllabus:" content (cells already existed)
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Value = "1. Basic notions of Marketing, Finance and Human Resources.2. Different organization settings."
$ws.Range("C16").Value = "1. Basic notions of Marketing, Finance and Human Resources.2. Different organization settings."

# Row 17 - now "Avaliação:" only (B/C cleared), reverts to default row height
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("B17:C17").Clear()
$ws.Rows(17).AutoFit()

# Row 18 - now "Método:" with professor placeholder text, becomes 60pt row.
# B18/C18 are brand-new cells for this row.
$ws.Range("A18").Value = "Método:"
$ws.Range("B10").Copy()
$ws.Range("B18").PasteSpecial($xlPasteFormats)
$ws.Range("B18").Value = "11079086 - Herlandí de Souza Andrade"
$ws.Range("C10").Copy()
$ws.Range("C18").PasteSpecial($xlPasteFormats)
$ws.Range("C18").Value = "11079086 - Herlandí de Souza Andrade"
$ws.Rows(18).RowHeight = 60

# Row 19 - label only changes to "Critério:"
$ws.Range("A19").Value = "Critério:"

# Row 20 - label only changes to "Norma de recuperação:"
$ws.Range("A20").Value = "Norma de recuperação:"

# Row 21 - label changes to "Bibliografia:", becomes 120pt row
$ws.Range("A21").Value = "Bibliografia:"
$ws.Rows(21).RowHeight = 120

# Row 22 (old long Bibliografia text row) is removed entirely
$ws.Rows(22).Delete()
